# Update countries & provincias Spain
# - Refresh case counts for a handful of countries
# - Re-sort the country table by "Casos totales" (column B) descending
# - Bump the "last updated" timestamp in the title cell

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Updated statistics for the affected countries --------------------
# Country name -> Casos totales, Nuevos casos, Casos activos, Recuperados,
#                 Casos criticos, Muertes hoy, Muertes
$updates = @{
    "Suiza"                  = @(11593, 696, 131, 11291, 141, 18, 171)
    "Noruega"                = @(3233, 149, 6, 3213, 70, 0, 14)
    "Israel"                 = @(2666, 297, 68, 2590, 39, 3, 8)
    "Dinamarca"              = @(1851, 127, 1, 1809, 94, 7, 41)
    "Finlandia"              = @(928, 48, 10, 913, 22, 2, 5)
    "India"                  = @(697, 40, 45, 638, 0, 2, 14)
    "Irak"                   = @(382, 36, 105, 241, 0, 7, 36)
    "Eslovaquia"             = @(226, 10, 2, 224, 2, 0, 0)
    "Republica de Macedonia" = @(201, 24, 3, 195, 1, 0, 3)
    "Malta"                  = @(134, 5, 2, 132, 1, 0, 0)
    "Georgia"                = @(79, 4, 10, 69, 1, 0, 0)
}

$searchRange = $ws.Range("A4:A204")

foreach ($country in $updates.Keys) {
    $vals = $updates[$country]
    $cell = $searchRange.Find($country)
    if ($cell -ne $null) {
        $row = $cell.Row
        $ws.Cells.Item($row, 2).Value = $vals[0]
        $ws.Cells.Item($row, 3).Value = $vals[1]
        $ws.Cells.Item($row, 4).Value = $vals[2]
        $ws.Cells.Item($row, 5).Value = $vals[3]
        $ws.Cells.Item($row, 6).Value = $vals[4]
        $ws.Cells.Item($row, 7).Value = $vals[5]
        $ws.Cells.Item($row, 8).Value = $vals[6]
    }
}

# --- 2. Re-sort the data rows (A4:H204) by Casos totales, descending -----
$dataRange = $ws.Range("A4:H204")
$sortKey = $ws.Range("B4:B204")
$dataRange.Sort($sortKey, 2, $null, $null, 1, $null, 1, 1)

# --- 3. Bump the "last updated" timestamp ---------------------------------
$ws.Range("A1").Value = "Datos actualizados a 26 de Marzo de 2020 a las 13:42"
